# Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts, etc.) across achievement / work
# experience bullet paragraphs, per the "Implement quantitative metrics
# highlighting across all resume formats" change.

$d = $word.ActiveDocument

# Word Range.Font.Color expects a BGR-packed long: 0x00BBGGRR.
# Target RGB 2C3E50 -> R=2C G=3E B=50 -> 0x00503E2C = 5258796
$metricColor = 5258796

function ApplyBoldColor($paraRange, $searchText) {
    if ($paraRange -eq $null) {
        Write-Host "WARNING: paragraph range is null, cannot search for '$searchText'"
        return $false
    }
    $sub = $paraRange.Duplicate
    $found = $sub.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $sub.Font.Bold = 1
        $sub.Font.Color = $metricColor
    } else {
        Write-Host "WARNING: could not find '$searchText' in target paragraph"
    }
    return $found
}

function FindParagraphRange($matchPattern) {
    $result = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $matchPattern) {
            $result = $p.Range
        }
    }
    if ($result -eq $null) {
        Write-Host "WARNING: no paragraph matched pattern '$matchPattern'"
    }
    return $result
}

# 1. Siege Analytics: "Discovered systematic race coding errors..." bullet
#    -> bold "23%" and "64%"
$p1 = FindParagraphRange "*Discovered systematic race coding errors*"
ApplyBoldColor $p1 "23%" | Out-Null
ApplyBoldColor $p1 "64%" | Out-Null

# 2. Siege Analytics: "Utilized advanced sampling methods..." bullet
#    -> bold "±4.2%", "±2.1%", "71%", "87%"
$p2 = FindParagraphRange "*Utilized advanced sampling methods to decrease survey margin of error*increasing voter turnout*"
ApplyBoldColor $p2 "±4.2%" | Out-Null
ApplyBoldColor $p2 "±2.1%" | Out-Null
ApplyBoldColor $p2 "71%" | Out-Null
ApplyBoldColor $p2 "87%" | Out-Null

# 3. "Trigonometric algorithm for boundary estimation..." bullet
#    -> bold "73.5%" and "$4.7M"
$p3 = FindParagraphRange "*Trigonometric algorithm for boundary estimation*"
ApplyBoldColor $p3 "73.5%" | Out-Null
ApplyBoldColor $p3 "$4.7M" | Out-Null

# 4. "Built real-time FEC analysis systems..." bullet -> bold "$2"
$p4 = FindParagraphRange "*Built real-time FEC analysis systems*"
ApplyBoldColor $p4 "$2" | Out-Null

# 5. Helm/Murmuration: "Modernized legacy ETL processes..." bullet -> bold "57%"
$p5 = FindParagraphRange "*Modernized legacy ETL processes*"
ApplyBoldColor $p5 "57%" | Out-Null

# 6. Key Achievements: "Predictive excellence..." bullet
#    -> bold "±4.2%" and "±2.1%"
$p6 = FindParagraphRange "*Predictive excellence*"
ApplyBoldColor $p6 "±4.2%" | Out-Null
ApplyBoldColor $p6 "±2.1%" | Out-Null

# 7. Key Achievements: "Increased voter turnout prediction accuracy..." bullet
#    -> bold "71%" and "87%"
$p7 = FindParagraphRange "*Increased voter turnout prediction accuracy*"
ApplyBoldColor $p7 "71%" | Out-Null
ApplyBoldColor $p7 "87%" | Out-Null

# 8. Key Achievements: "Methodological advancement..." bullet
#    -> bold "34%" and "28%"
$p8 = FindParagraphRange "*Methodological advancement*"
ApplyBoldColor $p8 "34%" | Out-Null
ApplyBoldColor $p8 "28%" | Out-Null
